$wb = $excel.ActiveWorkbook

# Add a new worksheet and name it "Birim"
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Birim"

# Move it to be the last sheet (after "Isler")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet) | Out-Null

# Re-fetch by name: the position-based reference may now point elsewhere after Move
$birim = $wb.Worksheets.Item("Birim")

# Fill in the data
$birim.Range("A1").Value = "Merkez"
$birim.Range("A2").Value = "İlçe"

# Leave selection on the cell right below the data, and make the new sheet
# the active / selected sheet (mirrors the other data sheets' state)
$birim.Range("A3").Select() | Out-Null
$birim.Select() | Out-Null
$birim.Activate() | Out-Null
